$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh of the cryptos price table: Price (D) and Volume(1h) (E) columns for
# rows 2-51, matching the latest coinranking.com snapshot pulled by the scheduled
# GitHub Actions job. All of these are plain text cells (not numeric), including
# the Price entries that happen to look like numbers (e.g. "246.77"), so for those
# we force a Text number format before writing the value - otherwise Excel would
# silently re-cast the text as a Number.

# --- Price column (D): values that look numeric -> force Text format first ---
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '246.77'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9997'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4729'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2927'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06517'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '22.48'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07772'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.7413'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '96.62'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '284.43'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.26'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007518'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.9992'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.316'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.9993'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.262'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.242'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '164.43'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.97'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.344'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.09768'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.300'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.191'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04905'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.135'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.707'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01899'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.842'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '76.17'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.281'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.016'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.4291'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9998'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8334'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '101.90'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.567'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.023'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '35.45'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '909.98'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05774'

# --- Remaining Price (D) and Volume(1h) (E) values: already unambiguous text ---
$ws.Range("D2").Value = '30.738.06'
$ws.Range("E2").Value = '  +0.84%  '
$ws.Range("D3").Value = '1.891.68'
$ws.Range("E3").Value = '  +1.18%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("E5").Value = '  -0.14%  '
$ws.Range("E6").Value = '  +0.02%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("E8").Value = '  +0.54%  '
$ws.Range("E9").Value = '  +0.54%  '
$ws.Range("E10").Value = '  +1.91%  '
$ws.Range("E11").Value = '  +0.70%  '
$ws.Range("E12").Value = '  -0.10%  '
$ws.Range("D13").Value = '1.885.81'
$ws.Range("E13").Value = '  +0.85%  '
$ws.Range("E14").Value = '  -0.65%  '
$ws.Range("E15").Value = '  +1.94%  '
$ws.Range("E16").Value = '  +3.98%  '
$ws.Range("D17").Value = '30.733.67'
$ws.Range("E17").Value = '  +0.95%  '
$ws.Range("E18").Value = '  -0.84%  '
$ws.Range("E19").Value = '  +0.18%  '
$ws.Range("E20").Value = '  -0.12%  '
$ws.Range("D21").Value = '2.134.55'
$ws.Range("E21").Value = '  +0.95%  '
$ws.Range("E22").Value = '  +1.66%  '
$ws.Range("E23").Value = '  +0.00%  '
$ws.Range("E24").Value = '  +1.67%  '
$ws.Range("E25").Value = '  -0.30%  '
$ws.Range("E26").Value = '  +0.64%  '
$ws.Range("E27").Value = '  +1.32%  '
$ws.Range("E28").Value = '  +0.31%  '
$ws.Range("E29").Value = '  -0.75%  '
$ws.Range("E30").Value = '  -2.14%  '
$ws.Range("E31").Value = '  -1.44%  '
$ws.Range("E32").Value = '  +0.80%  '
$ws.Range("E33").Value = '  +2.28%  '
$ws.Range("E34").Value = '  +2.20%  '
$ws.Range("E35").Value = '  +1.68%  '
$ws.Range("E36").Value = '  +0.81%  '
$ws.Range("E37").Value = '  -0.21%  '
$ws.Range("E38").Value = '  +2.77%  '
$ws.Range("E39").Value = '  +3.71%  '
$ws.Range("E40").Value = '  +5.00%  '
$ws.Range("E41").Value = '  +1.01%  '
$ws.Range("E42").Value = '  +2.58%  '
$ws.Range("E43").Value = '  +2.63%  '
$ws.Range("E44").Value = '  -0.03%  '
$ws.Range("E45").Value = '  +0.07%  '
$ws.Range("E46").Value = '  -0.17%  '
$ws.Range("E47").Value = '  +2.52%  '
$ws.Range("E48").Value = '  +1.00%  '
$ws.Range("E49").Value = '  +0.18%  '
$ws.Range("E50").Value = '  -1.23%  '
$ws.Range("E51").Value = '  +2.45%  '
